$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 8 (B8, C8, D8)
$ws.Range("B8").Value = 104
$ws.Range("C8").Value = 73
$ws.Range("D8").Value = 31

# Update existing values in row 9 (B9, C9, D9, E9)
$ws.Range("B9").Value = 120
$ws.Range("C9").Value = 45
$ws.Range("D9").Value = 58
$ws.Range("E9").Value = 18

# Delete column H entirely (H1:H10)
$ws.Range("H1:H10").Delete()

# Add new row 11 with decade bucket (2025, 2035]
$ws.Range("A11").Value = "(2025, 2035]"
$ws.Range("A1").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "(2025, 2035]"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
